# Daily attendance processing - 2025-10-07 22:17:25
# Reorders the "Recorded By" (column G) values for specific rows:
#   "System, backup@backdoor.com, system" -> "system, System, backup@backdoor.com"
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rowsSystemBackdoor = @(2, 29, 56)
foreach ($r in $rowsSystemBackdoor) {
    $ws.Cells.Item($r, 7).Value = "system, System, backup@backdoor.com"
}

$rowsSystemDnasr = @(3, 6, 11, 12, 13, 14, 15, 30, 33, 38, 39, 40, 41, 42, 57, 60, 65, 66, 67, 68, 69, 86, 89, 93, 112, 115, 119, 138, 141, 145)
foreach ($r in $rowsSystemDnasr) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}
